# Applies the "Added background to readme, calculated throughput for 10 repl
# system" commit to the workbook: extends the CvRDT-N (10-node replicated)
# throughput table with 5 extra zero rows + averages + an ops/sec throughput
# ratio, and adds the equivalent ops/sec throughput ratios (column P/R) to
# the MongoDB sheet, fixing up the CmRDT-O label indices as a side effect of
# dropping the now-unused "Plus 5 zeros" shared string.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "CvRDT-N": turn the "Plus 5 zeros" label (N23) into real zero data
# ---------------------------------------------------------------------
$wsN = $wb.Worksheets.Item("CvRDT-N")

# Row 23: replace the text label with numeric zeros across N:R
$wsN.Range("N23").Value = 0
$wsN.Range("O23").Value = 0
$wsN.Range("P23").Value = 0
$wsN.Range("Q23").Value = 0
$wsN.Range("R23").Value = 0

# Rows 24-27: four more rows of zeros (the "5 zeros" referred to by the old
# label: rows 23-27)
foreach ($r in 24..27) {
    $wsN.Range("N$r").Value = 0
    $wsN.Range("O$r").Value = 0
    $wsN.Range("P$r").Value = 0
    $wsN.Range("Q$r").Value = 0
    $wsN.Range("R$r").Value = 0
}

# Row 28: averages of the 10-row block N18:R27
$wsN.Range("N28").Formula = "=AVERAGE(N18:N27)"
$wsN.Range("O28").Formula = "=AVERAGE(O18:O27)"
$wsN.Range("P28").Formula = "=AVERAGE(P18:P27)"
$wsN.Range("Q28").Formula = "=AVERAGE(Q18:Q27)"
$wsN.Range("R28").Formula = "=AVERAGE(R18:R27)"

# Row 30: overall ops/sec throughput for the 10-replica system
$wsN.Range("O30").Formula = "=1543*100000"
$wsN.Range("O30").Font.Bold = $true
$wsN.Range("N30").Formula = "=O30/AVERAGE(N28:R28)"
$wsN.Range("N30").NumberFormat = "#,##0.00"

$wsN.Columns.Item(14).ColumnWidth = 10.28515625

$wsN.Range("O31").Select()

# ---------------------------------------------------------------------
# Sheet "MongoDB": derive ops/sec throughput ratios against the same
# 10-replica system throughput (R1)
# ---------------------------------------------------------------------
$wsM = $wb.Worksheets.Item("MongoDB")

$wsM.Range("R1").Formula = "=1543*100000"

$wsM.Range("P1").Formula = "=R1/B1"
$wsM.Range("P2").Formula = "=R1/E1"
$wsM.Range("P3").Formula = "=R1/H1"
$wsM.Range("P4").Formula = "=R1/K1"
$wsM.Range("P5").Formula = "=R1/N1"
$wsM.Range("P6").Formula = "=R1/B9"
$wsM.Range("P7").Formula = "=R1/E9"
$wsM.Range("P8").Formula = "=R1/H9"
$wsM.Range("P9").Formula = "=R1/K9"
$wsM.Range("P10").Formula = "=R1/N9"
$wsM.Range("P11").Formula = "=R1/B17"
$wsM.Range("P12").Formula = "=R1/E17"
$wsM.Range("P13").Formula = "=R1/H17"

$wsM.Columns.Item(16).ColumnWidth = 12.42578125
$wsM.Columns.Item(18).ColumnWidth = 10.28515625

$wsM.PageSetup.Orientation = 1

$wsM.Range("R1").Select()

$wsN.Activate()
